# Update the "dSF" column (F) values for rows 2-6 on the active sheet.
# This reflects a repull/recalculation of data (commit: "repull data, push all data, mean calculation").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 3
$ws.Range("F3").Value = 6
$ws.Range("F4").Value = 3
$ws.Range("F5").Value = -2
$ws.Range("F6").Value = -1
